$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet 1 (quality_comparison): build the two new border styles on C1/D1 ---
$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.Item(8).LineStyle = 1
$c1.Borders.Item(9).LineStyle = 1
$c1.Borders.Item(7).LineStyle = -4142
$c1.Borders.Item(10).LineStyle = -4142

$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.Item(8).LineStyle = 1
$d1.Borders.Item(7).LineStyle = -4142
$d1.Borders.Item(10).LineStyle = 1
$d1.Borders.Item(9).LineStyle = 1

# --- Sheet 2 (computational_comparison): reuse the same two border styles via
# copy/paste-special of formats, so the style table isn't bloated with
# duplicate/orphan xf entries ---
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)
$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$c1.Copy()
$ws2.Range("F1").PasteSpecial(-4122)
$d1.Copy()
$ws2.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Text / value edits (anonymise "fedcore" -> "approach", fix -0 -> 0) ---
$ws1.Range("C2").Value = "approach"
$ws1.Range("D5").Value = 0

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# G5 was an empty inline-string cell; drop it entirely (matches removal in diff)
$ws2.Range("G5").ClearContents()
